# Regen save_data to use K (strikeouts) instead of Strike# column values.
# Recalculated K values are written back into column G (rows 2-35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New strikeout (K) values, keyed by worksheet row number.
$kValues = [ordered]@{
    2  = 8
    3  = 6
    4  = 6
    5  = 4
    6  = 2
    7  = 4
    8  = 5
    9  = 7
    10 = 4
    11 = 6
    12 = 3
    13 = 4
    14 = 6
    15 = 4
    16 = 4
    17 = 7
    18 = 2
    19 = 2
    20 = 7
    21 = 6
    22 = 7
    23 = 5
    24 = 2
    25 = 7
    26 = 6
    27 = 7
    28 = 5
    29 = 1
    30 = 5
    31 = 3
    32 = 5
    33 = 8
    34 = 6
    35 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
